$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H4").Value = 538.0625
$ws.Range("I4").Value = 145.75
$ws.Range("J4").Value = 930.375
$ws.Range("K4").Value = 145.75
$ws.Range("L4").Value = 930.375
$ws.Range("M4").Value = -31.75
$ws.Range("N4").Value = -1158.375

$ws.Range("H5").Value = 41.384617
$ws.Range("I5").Value = 36.666668
$ws.Range("J5").Value = 45.42857
$ws.Range("K5").Value = 36.666668
$ws.Range("L5").Value = 45.42857
$ws.Range("M5").Value = 78.333332
$ws.Range("N5").Value = -275.42857

$ws.Range("H11").Value = 446809.97
$ws.Range("I11").Value = 446809.97
$ws.Range("K11").Value = 446809.97
$ws.Range("M11").Value = -446669.97

$ws.Range("H51").Value = 2741.6667
$ws.Range("J51").Value = 3557.1428
$ws.Range("L51").Value = 3557.1428
$ws.Range("N51").Value = -4525.1428

$ws.Range("H53").Value = 823.65216
$ws.Range("I53").Value = 1348.25
$ws.Range("K53").Value = 1348.25
$ws.Range("M53").Value = -711.25

$ws.Range("H70").Value = 2419.6
$ws.Range("J70").Value = 1774.5
$ws.Range("L70").Value = 5323.5
$ws.Range("N70").Value = -5863.5

$ws.Range("H73").Value = 2419.6
$ws.Range("J73").Value = 1774.5
$ws.Range("L73").Value = 5323.5
$ws.Range("N73").Value = -7195.5

$ws.Range("H103").Value = 1350
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1350
$ws.Range("K103").Value = 0
$ws.Range("N103").Value = -5222
$ws.Range("L103").Value = 4050
$ws.Range("M103").ClearContents()

$ws.Range("H116").Value = 8715.272000000001
$ws.Range("I116").Value = 8196.5
$ws.Range("K116").Value = 8196.5
$ws.Range("M116").Value = -4754.5

$ws.Range("H138").Value = 52686430
$ws.Range("J138").Value = 125003320
$ws.Range("L138").Value = 375009960
$ws.Range("N138").Value = -375020240

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 5447.4604
$ws.Range("I32").Value = 2674.9565
$ws.Range("K32").Value = 2674.9565
$ws.Range("M32").Value = -2387.9565

$ws.Range("H102").Value = 144639.58
$ws.Range("I102").Value = 168204.5
$ws.Range("K102").Value = 168204.5
$ws.Range("M102").Value = -166582.5

$ws.Range("H132").Value = 3753.2964
$ws.Range("I132").Value = 2910.7368
$ws.Range("K132").Value = 8732.2104
$ws.Range("M132").Value = -6202.2104

$ws.Range("H135").Value = 236357
$ws.Range("J135").Value = 236357
$ws.Range("L135").Value = 236357
$ws.Range("N135").Value = -246497

$ws.Range("H139").Value = 169886
$ws.Range("J139").Value = 169886
$ws.Range("L139").Value = 169886
$ws.Range("N139").Value = -180166

$ws.Range("H140").Value = 149999
$ws.Range("J140").Value = 149999
$ws.Range("L140").Value = 149999
$ws.Range("N140").Value = -160359

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H99").Value = 690921.5600000001
$ws.Range("I99").Value = 128331.125
$ws.Range("K99").Value = 128331.125
$ws.Range("M99").Value = -126833.125

$ws.Range("H107").Value = 4370.3335
$ws.Range("I107").Value = 5507.3335
$ws.Range("J107").Value = 3233.3333
$ws.Range("K107").Value = 5507.3335
$ws.Range("L107").Value = 3233.3333
$ws.Range("M107").Value = -3587.3335
$ws.Range("N107").Value = -7073.3333

$ws.Range("H134").Value = 3980.8572
$ws.Range("I134").Value = 2521.0908
$ws.Range("K134").Value = 7563.2724
$ws.Range("M134").Value = -5028.2724

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 9755.190000000001
$ws.Range("I7").Value = 9162.817999999999
$ws.Range("J7").Value = 10406.8
$ws.Range("K7").Value = 9162.817999999999
$ws.Range("L7").Value = 10406.8
$ws.Range("M7").Value = -9049.817999999999
$ws.Range("N7").Value = -10632.8

$ws.Range("H31").Value = 2729.8635
$ws.Range("I31").Value = 1911.9
$ws.Range("J31").Value = 3411.5
$ws.Range("K31").Value = 1911.9
$ws.Range("L31").Value = 3411.5
$ws.Range("M31").Value = -1616.9
$ws.Range("N31").Value = -4001.5

$ws.Range("H34").Value = 2729.8635
$ws.Range("I34").Value = 1911.9
$ws.Range("J34").Value = 3411.5
$ws.Range("K34").Value = 1911.9
$ws.Range("L34").Value = 3411.5
$ws.Range("M34").Value = -1709.9
$ws.Range("N34").Value = -3815.5

$ws.Range("H80").Value = 44088.8
$ws.Range("J80").Value = 44088.8
$ws.Range("L80").Value = 44088.8
$ws.Range("N80").Value = -46334.8

$ws.Range("H83").Value = 44088.8
$ws.Range("J83").Value = 44088.8
$ws.Range("L83").Value = 132266.4
$ws.Range("N83").Value = -143498.4

$ws.Range("H138").Value = 91304
$ws.Range("J138").Value = 91304
$ws.Range("L138").Value = 91304
$ws.Range("N138").Value = -101584

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H68").Value = 2794.7693
$ws.Range("I68").Value = 1001
$ws.Range("K68").Value = 3003
$ws.Range("M68").Value = -2192

$ws.Range("H71").Value = 2794.7693
$ws.Range("I71").Value = 1001
$ws.Range("K71").Value = 9009
$ws.Range("M71").Value = -4953

$ws.Range("H74").Value = 9687.4
$ws.Range("I74").Value = 219
$ws.Range("J74").Value = 15999.667
$ws.Range("K74").Value = 657
$ws.Range("L74").Value = 47999.001
$ws.Range("M74").Value = 404
$ws.Range("N74").Value = -50121.001

$ws.Range("H77").Value = 9687.4
$ws.Range("I77").Value = 219
$ws.Range("J77").Value = 15999.667
$ws.Range("K77").Value = 1971
$ws.Range("L77").Value = 143997.003
$ws.Range("M77").Value = 3333
$ws.Range("N77").Value = -154605.003

$ws.Range("H118").Value = 2311.5
$ws.Range("J118").Value = 3999.6667
$ws.Range("L118").Value = 11999.0001
$ws.Range("N118").Value = -14485.0001

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20576

$ws.Range("H43").Value = 12492.23
$ws.Range("J43").Value = 18750
$ws.Range("L43").Value = 18750
$ws.Range("N43").Value = -19052

$ws.Range("H46").Value = 33299.668
$ws.Range("I46").Value = 24949.5
$ws.Range("K46").Value = 24949.5
$ws.Range("M46").Value = -24793.5

$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996

$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984

$ws.Range("H126").Value = 3596.8125
$ws.Range("I126").Value = 2670.111
$ws.Range("K126").Value = 8010.333
$ws.Range("M126").Value = -5540.333

$ws.Range("H132").Value = 6683.1665
$ws.Range("I132").Value = 4869
$ws.Range("K132").Value = 14607
$ws.Range("M132").Value = -12077

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 85431.56
$ws.Range("I22").Value = 1299
$ws.Range("J22").Value = 192509.36
$ws.Range("K22").Value = 1299
$ws.Range("L22").Value = 192509.36
$ws.Range("M22").Value = -1004
$ws.Range("N22").Value = -193099.36

$ws.Range("H27").Value = 85431.56
$ws.Range("I27").Value = 1299
$ws.Range("J27").Value = 192509.36
$ws.Range("K27").Value = 1299
$ws.Range("L27").Value = 192509.36
$ws.Range("M27").Value = -1192
$ws.Range("N27").Value = -192723.36

$ws.Range("H40").Value = 4717462.5
$ws.Range("I40").Value = 104699.1
$ws.Range("J40").Value = 27781280
$ws.Range("K40").Value = 104699.1
$ws.Range("L40").Value = 27781280
$ws.Range("M40").Value = -104563.1
$ws.Range("N40").Value = -27781552

$ws.Range("H74").Value = 100197
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 100197
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H12").Value = 3007
$ws.Range("J12").Value = 3007
$ws.Range("L12").Value = 3007
$ws.Range("N12").Value = -3291

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H81").Value = 67906.75
$ws.Range("J81").Value = 132813
$ws.Range("L81").Value = 265626
$ws.Range("N81").Value = -267748

$ws.Range("H84").Value = 67906.75
$ws.Range("J84").Value = 132813
$ws.Range("L84").Value = 1328130
$ws.Range("N84").Value = -1338738

$ws.Range("H96").Value = 7536511.5
$ws.Range("I96").Value = 38667.668
$ws.Range("K96").Value = 38667.668
$ws.Range("M96").Value = -37294.668

$ws.Range("H122").Value = 3114.48
$ws.Range("I122").Value = 3156.15
$ws.Range("K122").Value = 9468.450000000001
$ws.Range("M122").Value = -7018.450000000001

$ws.Range("H126").Value = 3610.7646
$ws.Range("I126").Value = 3246.3
$ws.Range("J126").Value = 4131.4287
$ws.Range("K126").Value = 9738.900000000001
$ws.Range("L126").Value = 12394.2861
$ws.Range("M126").Value = -7268.900000000001
$ws.Range("N126").Value = -17334.2861

$ws.Range("H132").Value = 2677.2
$ws.Range("I132").Value = 2353.4285
$ws.Range("K132").Value = 7060.2855
$ws.Range("M132").Value = -4530.2855
